$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 362, which pushes the existing rows 362..487 down to 363..488.
$ws.Rows.Item(362).Insert()

# Populate the newly inserted row 362 with the new weekly record.
$ws.Cells.Item(362, 1).Value  = 9
$ws.Cells.Item(362, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(362, 3).Value  = "Metropolitana"
$ws.Cells.Item(362, 4).Value  = 44809
$ws.Cells.Item(362, 5).Value  = 13
$ws.Cells.Item(362, 6).Value  = 100112012
$ws.Cells.Item(362, 7).Value  = "Espinaca"
$ws.Cells.Item(362, 8).Value  = "Sin especificar"
$ws.Cells.Item(362, 9).Value  = "Primera"
$ws.Cells.Item(362, 10).Value = 70
$ws.Cells.Item(362, 11).Value = 6000
$ws.Cells.Item(362, 12).Value = 7000
$ws.Cells.Item(362, 13).Value = 6500
$ws.Cells.Item(362, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(362, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(362, 16).Value = 650
$ws.Cells.Item(362, 17).Value = 10
$ws.Cells.Item(362, 18).Value = "Hortaliza"
